$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7692.4546
$ws.Range("J40").Value = 4876.6
$ws.Range("L40").Value = 4876.6
$ws.Range("N40").Value = -5226.6
$ws.Range("H100").Value = 8275.764999999999
$ws.Range("I100").Value = 6409.8887
$ws.Range("J100").Value = 10374.875
$ws.Range("K100").Value = 6409.8887
$ws.Range("L100").Value = 10374.875
$ws.Range("M100").Value = -5868.8887
$ws.Range("N100").Value = -11456.875
$ws.Range("H112").Value = 1657.3158
$ws.Range("J112").Value = 1815.1428
$ws.Range("L112").Value = 5445.428400000001
$ws.Range("N112").Value = -7661.428400000001
$ws.Range("H132").Value = 10143.862
$ws.Range("I132").Value = 1351.909
$ws.Range("K132").Value = 4055.727
$ws.Range("M132").Value = -1525.727
$ws.Range("H138").Value = 4638.643
$ws.Range("J138").Value = 4891.0557
$ws.Range("L138").Value = 14673.1671
$ws.Range("N138").Value = -24953.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14929825
$ws.Range("I32").Value = 16398195
$ws.Range("K32").Value = 16398195
$ws.Range("M32").Value = -16397908
$ws.Range("H74").Value = 3987
$ws.Range("I74").Value = 3324.25
$ws.Range("J74").Value = 5501.857
$ws.Range("K74").Value = 3324.25
$ws.Range("L74").Value = 5501.857
$ws.Range("M74").Value = -2450.25
$ws.Range("N74").Value = -7249.857
$ws.Range("H77").Value = 3987
$ws.Range("I77").Value = 3324.25
$ws.Range("J77").Value = 5501.857
$ws.Range("K77").Value = 16621.25
$ws.Range("L77").Value = 27509.285
$ws.Range("M77").Value = -12253.25
$ws.Range("N77").Value = -36245.285
$ws.Range("H132").Value = 2748.25
$ws.Range("I132").Value = 2011.5
$ws.Range("K132").Value = 6034.5
$ws.Range("M132").Value = -3504.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7915.857
$ws.Range("I107").Value = 8082.2
$ws.Range("K107").Value = 8082.2
$ws.Range("M107").Value = -6162.2
$ws.Range("H134").Value = 7759.4517
$ws.Range("I134").Value = 2546.074
$ws.Range("K134").Value = 7638.222
$ws.Range("M134").Value = -5103.222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 7725
$ws.Range("I21").Value = 900
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 900
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -665
$ws.Range("N21").Value = -10470
$ws.Range("H58").Value = 2644.7646
$ws.Range("I58").Value = 2264.889
$ws.Range("J58").Value = 3072.125
$ws.Range("K58").Value = 2264.889
$ws.Range("L58").Value = 3072.125
$ws.Range("M58").Value = -2061.889
$ws.Range("N58").Value = -3478.125
$ws.Range("H134").Value = 2262.1143
$ws.Range("I134").Value = 2262.1143
$ws.Range("K134").Value = 6786.342900000001
$ws.Range("M134").Value = -4251.342900000001
$ws.Range("H136").Value = 2644.7646
$ws.Range("I136").Value = 2264.889
$ws.Range("J136").Value = 3072.125
$ws.Range("K136").Value = 6794.667
$ws.Range("L136").Value = 9216.375
$ws.Range("M136").Value = -4244.667
$ws.Range("N136").Value = -14316.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 32500
$ws.Range("J57").Value = 32500
$ws.Range("L57").Value = 32500
$ws.Range("N57").Value = -34140
$ws.Range("H103").Value = 19977.223
$ws.Range("J103").Value = 19977.223
$ws.Range("L103").Value = 19977.223
$ws.Range("N103").Value = -22321.223
$ws.Range("H128").Value = 49769.23
$ws.Range("J128").Value = 49769.23
$ws.Range("L128").Value = 49769.23
$ws.Range("N128").Value = -59729.23
$ws.Range("H132").Value = 7453.25
$ws.Range("I132").Value = 7135.0386
$ws.Range("K132").Value = 21405.1158
$ws.Range("M132").Value = -18875.1158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5693.25
$ws.Range("J7").Value = 5082.3335
$ws.Range("L7").Value = 5082.3335
$ws.Range("N7").Value = -5306.3335
$ws.Range("H16").Value = 1272
$ws.Range("I16").Value = 1197.0834
$ws.Range("J16").Value = 1721.5
$ws.Range("K16").Value = 1197.0834
$ws.Range("L16").Value = 1721.5
$ws.Range("M16").Value = -1027.0834
$ws.Range("N16").Value = -2061.5
$ws.Range("H40").Value = 3894
$ws.Range("I40").Value = 3948.889
$ws.Range("J40").Value = 3400
$ws.Range("K40").Value = 3948.889
$ws.Range("L40").Value = 3400
$ws.Range("M40").Value = -3812.889
$ws.Range("N40").Value = -3672
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -798
$ws.Range("N61").ClearContents()
$ws.Range("H68").Value = 2878.8
$ws.Range("I68").Value = 2878.8
$ws.Range("K68").Value = 2878.8
$ws.Range("M68").Value = -2129.8
$ws.Range("H71").Value = 2878.8
$ws.Range("I71").Value = 2878.8
$ws.Range("K71").Value = 14394
$ws.Range("M71").Value = -10650
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("H93").Value = 1424.6818
$ws.Range("I93").Value = 1418.9474
$ws.Range("K93").Value = 1418.9474
$ws.Range("M93").Value = -170.9474
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1170
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 4699.2144
$ws.Range("I122").Value = 4299.091
$ws.Range("K122").Value = 12897.273
$ws.Range("M122").Value = -10447.273
$ws.Range("H126").Value = 5693.25
$ws.Range("J126").Value = 5082.3335
$ws.Range("L126").Value = 15247.0005
$ws.Range("N126").Value = -20187.0005
$ws.Range("H128").Value = 50000
$ws.Range("J128").Value = 50000
$ws.Range("L128").Value = 50000
$ws.Range("N128").Value = -59960
$ws.Range("H132").Value = 4693.617
$ws.Range("I132").Value = 4855.5366
$ws.Range("K132").Value = 14566.6098
$ws.Range("M132").Value = -12036.6098

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()
$ws.Range("H110").Value = 57496
$ws.Range("J110").Value = 57496
$ws.Range("L110").Value = 57496
$ws.Range("N110").Value = -65676
$ws.Range("H120").Value = 77999.5
$ws.Range("J120").Value = 77999.5
$ws.Range("L120").Value = 77999.5
$ws.Range("N120").Value = -87675.5
$ws.Range("H132").Value = 4333.8086
$ws.Range("I132").Value = 3795.543
$ws.Range("J132").Value = 5903.75
$ws.Range("K132").Value = 11386.629
$ws.Range("L132").Value = 17711.25
$ws.Range("M132").Value = -8856.629000000001
$ws.Range("N132").Value = -22771.25
$ws.Range("H136").Value = 1566852.6
$ws.Range("I136").Value = 2383912.8
$ws.Range("K136").Value = 7151738.399999999
$ws.Range("M136").Value = -7149188.399999999
$ws.Range("H140").Value = 59998.75
$ws.Range("J140").Value = 59998.75
$ws.Range("L140").Value = 59998.75
$ws.Range("N140").Value = -70358.75
